$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-18 12:18:26"
$ws.Range("E3").Value = "2026-02-18 12:18:28"
$ws.Range("E4").Value = "2026-02-18 12:18:30"
$ws.Range("E5").Value = "2026-02-18 12:18:33"
$ws.Range("E6").Value = "2026-02-18 12:18:35"
$ws.Range("E7").Value = "2026-02-18 12:18:38"
$ws.Range("E8").Value = "2026-02-18 12:18:40"
$ws.Range("E9").Value = "2026-02-18 12:18:42"
$ws.Range("E10").Value = "2026-02-18 12:18:45"
$ws.Range("E11").Value = "2026-02-18 12:18:47"
$ws.Range("E12").Value = "2026-02-18 12:18:49"
$ws.Range("E13").Value = "2026-02-18 12:18:52"
$ws.Range("E14").Value = "2026-02-18 12:18:54"
$ws.Range("E15").Value = "2026-02-18 12:18:56"
$ws.Range("E16").Value = "2026-02-18 12:18:58"
$ws.Range("E17").Value = "2026-02-18 12:19:01"
$ws.Range("E18").Value = "2026-02-18 12:19:03"
$ws.Range("E19").Value = "2026-02-18 12:19:06"
$ws.Range("E20").Value = "2026-02-18 12:19:08"
$ws.Range("E21").Value = "2026-02-18 12:19:11"
$ws.Range("E22").Value = "2026-02-18 12:19:13"
$ws.Range("E23").Value = "2026-02-18 12:19:15"
$ws.Range("E24").Value = "2026-02-18 12:19:18"
$ws.Range("E25").Value = "2026-02-18 12:19:20"
$ws.Range("E26").Value = "2026-02-18 12:19:22"
$ws.Range("E27").Value = "2026-02-18 12:19:24"
$ws.Range("E28").Value = "2026-02-18 12:19:27"
$ws.Range("E29").Value = "2026-02-18 12:19:29"
$ws.Range("E30").Value = "2026-02-18 12:19:31"
$ws.Range("E31").Value = "2026-02-18 12:19:33"
$ws.Range("E32").Value = "2026-02-18 12:19:36"
$ws.Range("E33").Value = "2026-02-18 12:19:38"
$ws.Range("E34").Value = "2026-02-18 12:19:41"
$ws.Range("E35").Value = "2026-02-18 12:19:43"
$ws.Range("E36").Value = "2026-02-18 12:19:45"
$ws.Range("E37").Value = "2026-02-18 12:19:48"
$ws.Range("E38").Value = "2026-02-18 12:19:50"
$ws.Range("E39").Value = "2026-02-18 12:19:52"
$ws.Range("E40").Value = "2026-02-18 12:19:55"
$ws.Range("E41").Value = "2026-02-18 12:19:57"
$ws.Range("E42").Value = "2026-02-18 12:19:59"
$ws.Range("E43").Value = "2026-02-18 12:20:02"
$ws.Range("E44").Value = "2026-02-18 12:20:04"
$ws.Range("E45").Value = "2026-02-18 12:20:07"
$ws.Range("E46").Value = "2026-02-18 12:20:09"
